$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target state (row => Coin, Link, Price, Volume(1h)) for rows 2-51 after the data refresh.
$data = @(
    @{Row=2; B="Bitcoin"; C="https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"; D="76.024.05"; E="  +1.56%  "},
    @{Row=3; B="Ethereum"; C="https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"; D="2.909.69"; E="  +3.42%  "},
    @{Row=4; B="TetherUSD"; C="https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"; D="1.00"; E="  +0.01%  "},
    @{Row=5; B="Solana"; C="https://coinranking.com/coin/zNZHO_Sjf+solana-sol"; D="198.66"; E="  +5.58%  "},
    @{Row=6; B="BNB"; C="https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"; D="600.07"; E="  +0.57%  "},
    @{Row=7; B="USDC"; C="https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"; D="0.999"; E="  -0.07%  "},
    @{Row=8; B="XRP"; C="https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"; D="0.549"; E="  -1.33%  "},
    @{Row=9; B="Dogecoin"; C="https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"; D="0.200"; E="  +4.12%  "},
    @{Row=10; B="LidoStakedEther"; C="https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"; D="2.908.29"; E="  +3.41%  "},
    @{Row=11; B="Cardano"; C="https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"; D="0.430"; E="  +16.53%  "},
    @{Row=12; B="TRON"; C="https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"; D="0.161"; E="  -0.17%  "},
    @{Row=13; B="Toncoin"; C="https://coinranking.com/coin/67YlI0K1b+toncoin-ton"; D="4.87"; E="  +0.94%  "},
    @{Row=14; B="WrappedliquidstakedEther2.0"; C="https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"; D="3.443.17"; E="  +3.49%  "},
    @{Row=15; B="WrappedBTC"; C="https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"; D="75.800.82"; E="  +1.31%  "},
    @{Row=16; B="ShibaInu"; C="https://coinranking.com/coin/xz24e0BjL+shibainu-shib"; D="0.0000193"; E="  +3.52%  "},
    @{Row=17; B="Avalanche"; C="https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"; D="27.39"; E="  +1.50%  "},
    @{Row=18; B="WrappedEther"; C="https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"; D="2.907.44"; E="  +3.41%  "},
    @{Row=19; B="Chainlink"; C="https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"; D="12.97"; E="  +5.48%  "},
    @{Row=20; B="Uniswap"; C="https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"; D="8.75"; E="  -2.04%  "},
    @{Row=21; B="BitcoinCash"; C="https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"; D="370.22"; E="  -1.16%  "},
    @{Row=22; B="SuiNetwork"; C="https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"; D="2.30"; E="  +2.61%  "},
    @{Row=23; B="Polkadot"; C="https://coinranking.com/coin/25W7FG7om+polkadot-dot"; D="4.29"; E="  +5.09%  "},
    @{Row=24; B="Dai"; C="https://coinranking.com/coin/MoTuySvg7+dai-dai"; D="1.00"; E="  -0.13%  "},
    @{Row=25; B="Litecoin"; C="https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"; D="71.17"; E="  +0.72%  "},
    @{Row=26; B="WrappedeETH"; C="https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"; D="3.059.92"; E="  +3.42%  "},
    @{Row=27; B="NEARProtocol"; C="https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"; D="4.20"; E="  +1.38%  "},
    @{Row=28; B="Aptos"; C="https://coinranking.com/coin/HGYj5JCv5+aptos-apt"; D="9.65"; E="  +1.11%  "},
    @{Row=29; B="PEPE"; C="https://coinranking.com/coin/03WI8NQPF+pepe-pepe"; D="0.0000109"; E="  +6.12%  "},
    @{Row=30; B="Binance-PegBSC-USD"; C="https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"; D="0.996"; E="  -0.27%  "},
    @{Row=31; B="Fetch.AI"; C="https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"; D="1.41"; E="  +1.84%  "},
    @{Row=32; B="Bittensor"; C="https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"; D="502.26"; E="  -2.18%  "},
    @{Row=33; B="InternetComputer(DFINITY)"; C="https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"; D="7.70"; E="  -1.95%  "},
    @{Row=34; B="PancakeSwap"; C="https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"; D="1.82"; E="  +1.30%  "},
    @{Row=35; B="FirstDigitalUSD"; C="https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"; D="0.999"; E="  -0.04%  "},
    @{Row=36; B="Monero"; C="https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"; D="165.10"; E="  +1.34%  "},
    @{Row=37; B="EthereumClassic"; C="https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"; D="20.22"; E="  +0.82%  "},
    @{Row=38; B="WhiteBITCoin"; C="https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"; D="19.64"; E="  +1.66%  "},
    @{Row=39; B="Cronos"; C="https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"; D="0.106"; E="  +22.88%  "},
    @{Row=40; B="Kaspa"; C="https://coinranking.com/coin/V8GxkwWow+kaspa-kas"; D="0.113"; E="  -4.68%  "},
    @{Row=41; B="USDe"; C="https://coinranking.com/coin/exbfr2U-0+usde-usde"; D="1.00"; E="  -0.07%  "},
    @{Row=42; B="Aave"; C="https://coinranking.com/coin/ixgUfzmLR+aave-aave"; D="180.47"; E="  -0.95%  "},
    @{Row=43; B="PolygonEcosystemToken"; C="https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"; D="0.346"; E="  +2.50%  "},
    @{Row=44; B="RenderToken"; C="https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"; D="4.99"; E="  -0.72%  "},
    @{Row=45; B="Stacks"; C="https://coinranking.com/coin/mMPrMcB7+stacks-stx"; D="1.65"; E="  -1.47%  "},
    @{Row=46; B="OKB"; C="https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"; D="40.06"; E="  +0.59%  "},
    @{Row=47; B="ImmutableX"; C="https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; D="1.19"; E="  -1.78%  "},
    @{Row=48; B="dogwifhat"; C="https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"; D="2.33"; E="  -0.97%  "},
    @{Row=49; B="ARBITRUM"; C="https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"; D="0.572"; E="  +1.06%  "},
    @{Row=50; B="Filecoin"; C="https://coinranking.com/coin/ymQub4fuB+filecoin-fil"; D="3.72"; E="  -0.04%  "},
    @{Row=51; B="Mantle"; C="https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"; D="0.657"; E="  +6.59%  "}
)

foreach ($row in $data) {
    $r = $row.Row

    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C

    # Column D (Price) holds values like "1.00" or "76.024.05" that Excel would
    # otherwise reinterpret as numbers, so force text storage, then drop back to
    # the default (unstyled) cell style once the literal text value is set.
    $dCell = $ws.Cells.Item($r, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $row.D
    $dCell.Style = "Normal"

    $ws.Cells.Item($r, 5).Value = $row.E
}
